# Adds the "Email" attribute to the Utente entity in the E-R schema.
# - Moves the existing "Data Nascita" attribute marker + label slightly
#   to make room for the new attribute.
# - Duplicates the attribute-marker group (connector + oval "crow's foot")
#   and the attribute-label textbox, repositions them, and relabels the
#   textbox "Email".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- locate the existing shapes by their (stable) shape Id -----------------
$markerGroup = $null
$dataNascitaLabel = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 309) { $markerGroup = $shp }
    if ($shp.Id -eq 323) { $dataNascitaLabel = $shp }
}

# --- 1. nudge the "Data Nascita" marker + label to their new spot ----------
$markerGroup.Left = 589.9580078125
$markerGroup.Top = 325.9753723144531

$dataNascitaLabel.Left = 553.7123413085938
$dataNascitaLabel.Top = 345.430419921875

# --- 2. duplicate the marker group for the new "Email" attribute -----------
$newMarkerRange = $markerGroup.Duplicate()
$newMarkerGroup = $newMarkerRange.Item(1)
$newMarkerGroup.Name = "Gruppo 111"
$newMarkerGroup.Left = 612.8131496062992
$newMarkerGroup.Top = 330.5576377952756
$newMarkerGroup.Rotation = 206.14831666666666

$newMarkerGroup.GroupItems.Item(1).Name = "Connettore 1 112"
$newMarkerGroup.GroupItems.Item(2).Name = "Ovale 114"

# --- 3. duplicate the label textbox and turn it into "Email" ---------------
$newLabelRange = $dataNascitaLabel.Duplicate()
$newLabel = $newLabelRange.Item(1)
$newLabel.Name = "CasellaDiTesto 118"
$newLabel.Left = 592.7207641601562
$newLabel.Top = 354.13189697265625
$newLabel.Width = 39.467403411865234
$newLabel.Height = 14.540630340576172
$newLabel.TextFrame.TextRange.Text = "Email"
